$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (e.g. "7.86", "3.309.94")
# that must stay plain text, matching the source inlineStr cells. Force the
# cell number format to Text before assignment so Excel does not auto-coerce
# them into numbers (which would also mangle multi-dot values like "60.300.97").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.300.97"
$ws.Range("E2").Value = "  -2.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.309.94"
$ws.Range("E3").Value = "  -2.75%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.12"
$ws.Range("E5").Value = "  -3.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.41"
$ws.Range("E6").Value = "  -3.51%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.310.45"
$ws.Range("E8").Value = "  -2.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.86"
$ws.Range("E10").Value = "  -1.47%  "

$ws.Range("E11").Value = "  -3.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.882.34"
$ws.Range("E13").Value = "  -2.75%  "

$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.97"
$ws.Range("E15").Value = "  -4.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.312.50"
$ws.Range("E16").Value = "  -2.76%  "

$ws.Range("E17").Value = "  -2.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.319.02"
$ws.Range("E18").Value = "  -2.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("E19").Value = "  -2.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.45"
$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.66"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.93"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.14"
$ws.Range("E23").Value = "  -0.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.544"
$ws.Range("E24").Value = "  -3.43%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.441.59"
$ws.Range("E26").Value = "  -3.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000103"
$ws.Range("E27").Value = "  -6.95%  "

$ws.Range("E28").Value = "  -4.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").Value = "  -4.65%  "

$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("E33").Value = "  -3.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.57"
$ws.Range("E34").Value = "  -2.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -4.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("E36").Value = "  -4.62%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.58"
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("E38").Value = "  -5.50%  "

$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.04"
$ws.Range("E40").Value = "  -12.97%  "

$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.341.32"
$ws.Range("E41").Value = "  -2.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0739"
$ws.Range("E42").Value = "  -5.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.01"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.752"
$ws.Range("E44").Value = "  -3.79%  "

$ws.Range("E45").Value = "  -3.32%  "

$ws.Range("E46").Value = "  -4.07%  "

$ws.Range("E47").Value = "  -3.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.369.10"
$ws.Range("E48").Value = "  -6.72%  "

$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("E50").Value = "  -5.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.33"
$ws.Range("E51").Value = "  -4.96%  "
